$d = $word.ActiveDocument

# The subtitle paragraph currently reads "Modest Software Engineering Project "
# split across three runs: "Modest ", "Software", " Engineering Project ".
# We need to remove the leading "Modest " run entirely, leaving
# "Software Engineering Project " (still as two separate runs).
#
# A plain Find/Delete over the "Modest " text causes this runtime to
# auto-merge the now-adjacent "Software" and " Engineering Project " runs
# (they share identical run formatting). To keep them distinct - matching
# the original document structure - we briefly perturb the trailing run's
# font size so it no longer matches its neighbour's formatting while the
# "Modest " run is removed, then restore the original size afterwards.

$tail = $d.Content
$tail.Find.Execute(" Engineering Project ")
$originalSize = $tail.Font.Size
$tail.Font.Size = $originalSize + 1

$lead = $d.Content
$lead.Find.Execute("Modest ")
$lead.Delete()

$tail2 = $d.Content
$tail2.Find.Execute(" Engineering Project ")
$tail2.Font.Size = $originalSize
